# Scheduled runner update: refresh market-derived figures (columns H-N)
# on the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR "Profits" sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 10168
$ws.Range("I53").Value = 17160.572
$ws.Range("K53").Value = 17160.572
$ws.Range("M53").Value = -16523.572
$ws.Range("H58").Value = 1012
$ws.Range("I58").Value = 375.27274
$ws.Range("J58").Value = 3346.6667
$ws.Range("K58").Value = 1125.81822
$ws.Range("L58").Value = 10040.0001
$ws.Range("M58").Value = -975.8182200000001
$ws.Range("N58").Value = -10340.0001
$ws.Range("H112").Value = 4997.3076
$ws.Range("I112").Value = 1200
$ws.Range("J112").Value = 5687.727
$ws.Range("K112").Value = 3600
$ws.Range("L112").Value = 17063.181
$ws.Range("M112").Value = -2492
$ws.Range("N112").Value = -19279.181
$ws.Range("H135").Value = 810.2222
$ws.Range("I135").Value = 724
$ws.Range("K135").Value = 6516
$ws.Range("M135").Value = -3981
$ws.Range("H138").Value = 2881.125
$ws.Range("I138").Value = 2959.3809
$ws.Range("K138").Value = 8878.1427
$ws.Range("M138").Value = -3738.1427

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6379.1875
$ws.Range("I32").Value = 5024.263
$ws.Range("K32").Value = 5024.263
$ws.Range("M32").Value = -4737.263
$ws.Range("H97").Value = 996.2
$ws.Range("I97").Value = 495.25
$ws.Range("J97").Value = 3000
$ws.Range("K97").Value = 495.25
$ws.Range("L97").Value = 3000
$ws.Range("M97").Value = 0.75
$ws.Range("N97").Value = -3992

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H134").Value = 7239.9536
$ws.Range("I134").Value = 7841
$ws.Range("K134").Value = 23523
$ws.Range("M134").Value = -20988
$ws.Range("H140").Value = 50000
$ws.Range("J140").Value = 50000
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 999.5
$ws.Range("I105").Value = 999.5
$ws.Range("K105").Value = 999.5
$ws.Range("M105").Value = 747.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 720.8
$ws.Range("I5").Value = 614.7143
$ws.Range("K5").Value = 1844.1429
$ws.Range("M5").Value = -1732.1429
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("N11").ClearContents()
$ws.Range("H22").Value = 2667
$ws.Range("J22").Value = 4000
$ws.Range("L22").Value = 12000
$ws.Range("N22").Value = -12338
$ws.Range("H27").Value = 2667
$ws.Range("J27").Value = 4000
$ws.Range("L27").Value = 12000
$ws.Range("N27").Value = -12204
$ws.Range("H45").Value = 979.75
$ws.Range("J45").Value = 979.75
$ws.Range("L45").Value = 2939.25
$ws.Range("N45").Value = -4003.25
$ws.Range("H122").Value = 879.6667
$ws.Range("I122").Value = 547.2
$ws.Range("K122").Value = 4924.8
$ws.Range("M122").Value = -2474.8
$ws.Range("H128").Value = 400000
$ws.Range("I128").Value = 400000
$ws.Range("K128").Value = 1200000
$ws.Range("M128").Value = -1195020
$ws.Range("H131").Value = 760.5979599999999
$ws.Range("J131").Value = 804.3977
$ws.Range("L131").Value = 2413.1931
$ws.Range("N131").Value = -12493.1931
$ws.Range("H135").Value = 720.8
$ws.Range("I135").Value = 614.7143
$ws.Range("K135").Value = 5532.428699999999
$ws.Range("M135").Value = -2997.428699999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2199.8572
$ws.Range("I80").Value = 2299.5
$ws.Range("J80").Value = 2160
$ws.Range("K80").Value = 2299.5
$ws.Range("L80").Value = 2160
$ws.Range("M80").Value = -1301.5
$ws.Range("N80").Value = -4156
$ws.Range("H83").Value = 2199.8572
$ws.Range("I83").Value = 2299.5
$ws.Range("J83").Value = 2160
$ws.Range("K83").Value = 11497.5
$ws.Range("L83").Value = 10800
$ws.Range("M83").Value = -6505.5
$ws.Range("N83").Value = -20784

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 10676
$ws.Range("J43").Value = 10676
$ws.Range("L43").Value = 10676
$ws.Range("N43").Value = -11062
$ws.Range("H75").Value = 15000
$ws.Range("I75").Value = 15000
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 15000
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("M75").Value = -14064
$ws.Range("H78").Value = 15000
$ws.Range("I78").Value = 15000
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 45000
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("M78").Value = -40320
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H132").Value = 2338.9092
$ws.Range("I132").Value = 1946
$ws.Range("K132").Value = 5838
$ws.Range("M132").Value = -3308
$ws.Range("H136").Value = 2967.2222
$ws.Range("I136").Value = 1692.1875
$ws.Range("K136").Value = 5076.5625
$ws.Range("M136").Value = -2526.5625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("H103").Value = 19999.5
$ws.Range("J103").Value = 19999.5
$ws.Range("L103").Value = 19999.5
$ws.Range("N103").Value = -22343.5
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
